# Apply updated crypto price / volume(1h) figures per the Dec 9 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.293.68'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '2.365.84'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'0.694"
$ws.Range('E5').Value = '  +5.78%  '
$ws.Range('D6').Value = "'243.32"
$ws.Range('E6').Value = '  +3.25%  '
$ws.Range('D7').Value = "'74.31"
$ws.Range('E7').Value = '  +3.00%  '
$ws.Range('D9').Value = "'0.601"
$ws.Range('E9').Value = '  +28.02%  '
$ws.Range('D10').Value = "'0.104"
$ws.Range('E10').Value = '  +6.35%  '
$ws.Range('E11').Value = '  +16.69%  '
$ws.Range('E12').Value = '  +20.25%  '
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').Value = '2.719.12'
$ws.Range('E14').Value = '  -0.78%  '
$ws.Range('D15').Value = "'17.03"
$ws.Range('E15').Value = '  +7.56%  '
$ws.Range('D16').Value = "'0.919"
$ws.Range('E16').Value = '  +7.30%  '
$ws.Range('D17').Value = '2.355.41'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').Value = '44.257.28'
$ws.Range('E18').Value = '  +1.70%  '
$ws.Range('D19').Value = "'0.0000105"
$ws.Range('E19').Value = '  +5.15%  '
$ws.Range('D20').Value = "'6.74"
$ws.Range('E20').Value = '  +5.62%  '
$ws.Range('D21').Value = "'78.67"
$ws.Range('E21').Value = '  +5.33%  '
$ws.Range('D22').Value = "'257.30"
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = "'2.57"
$ws.Range('E24').Value = '  +3.56%  '
$ws.Range('E25').Value = '  -3.11%  '
$ws.Range('D26').Value = "'10.82"
$ws.Range('E26').Value = '  +7.88%  '
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = "'1.64"
$ws.Range('E28').Value = '  +7.08%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = "'22.73"
$ws.Range('E29').Value = '  -0.63%  '
$ws.Range('D30').Value = "'175.26"
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('E31').Value = '  +1.90%  '
$ws.Range('D32').Value = "'0.136"
$ws.Range('E32').Value = '  +6.39%  '
$ws.Range('E33').Value = '  +8.29%  '
$ws.Range('E34').Value = '  +10.17%  '
$ws.Range('E35').Value = '  +6.66%  '
$ws.Range('D36').Value = "'3.91"
$ws.Range('E36').Value = '  +5.42%  '
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('D38').Value = "'6.56"
$ws.Range('E38').Value = '  -0.83%  '
$ws.Range('D39').Value = "'0.0276"
$ws.Range('E39').Value = '  +7.71%  '
$ws.Range('D40').Value = "'19.22"
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').Value = "'9.12"
$ws.Range('E41').Value = '  +2.58%  '
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').Value = "'0.199"
$ws.Range('E43').Value = '  +17.90%  '
$ws.Range('E44').Value = '  +5.23%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').Value = "'1.26"
$ws.Range('E45').Value = '  +3.42%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = "'2.51"
$ws.Range('E46').Value = '  +12.21%  '
$ws.Range('D47').Value = "'101.36"
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('D50').Value = '1.465.06'
$ws.Range('E50').Value = '  +1.01%  '
$ws.Range('D51').Value = "'53.52"
$ws.Range('E51').Value = '  +5.32%  '
